$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 ("风" element), mirroring row 2 ("火") but left a blank row 3
# in between (robustness on blank lines).
$ws.Cells.Item(4, 1).Value = 114513
$ws.Cells.Item(4, 2).Value = "风"
$ws.Cells.Item(4, 3).Value = "test"
$ws.Cells.Item(4, 4).Value = "巫师"
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = "4火2水1风"
$ws.Cells.Item(4, 7).Value = "4水"
$ws.Cells.Item(4, 8).Value = "变成一坨大便"
$ws.Cells.Item(4, 9).Value = "牢大，我想你了\n牢大牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了牢大，我想你了"

# Give the whole new row an explicit black font (this is what produced the
# extra font/cellXf entries in styles.xml) and make sure J4 participates too
# even though it stays empty.
$ws.Range("A4:J4").Font.Color = 0

# Page setup tweak.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved after editing.
$ws.Range("G11").Select() | Out-Null
